# Remove the six DirectDebitPaymentRequest rows that were pulled out of the
# export (IDs 5234, 5235, 5236, 5244, 5253, 5254 -> worksheet rows
# 8, 9, 10, 13, 20, 21). Deleting from the bottom up keeps the remaining
# row numbers stable while each Delete() runs.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rowsToDelete = @(21, 20, 13, 10, 9, 8)
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}
